$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a run split at a given character offset by briefly adding and
# then removing a bookmark there. Word (and this host) never re-coalesces two
# runs that happened to be separated by a bookmark, even after the bookmark
# itself is removed, so this lets us reproduce an exact run boundary between
# two adjacent runs that otherwise share identical formatting.
# ---------------------------------------------------------------------------
function Split-RunAt($pos, $markName) {
    $bm = $d.Range($pos, $pos)
    $d.Bookmarks.Add($markName, $bm)
    $d.Bookmarks.Item($markName).Delete()
}

# ===========================================================================
# 1) Heading paragraph: "BACKGROUND AND GOALS FOR FELLOWSHIP TRAINING (six
#    page limit)" -> "RESPECTIVE CONTRIBUTIONS (one page limit)"
# ===========================================================================
$p1 = $d.Paragraphs.Item(1)

# Replace "...TRAINING " (including the trailing space that used to belong to
# the following " (" run) with "RESPECTIVE CONTRIBUTIONS " so only a single
# space remains before the parenthesis, matching the target text exactly.
$rngTitle = $d.Content
$rngTitle.Find.Execute("BACKGROUND AND GOALS FOR FELLOWSHIP TRAINING ", $false, $false, $false, $false, $false, $true, 1, $false, "RESPECTIVE CONTRIBUTIONS ", 2)

# Split "RESPECTIVE CONTRIBUTIONS " from the following "(" so they remain two
# distinct runs (matching the target markup).
$txt = $p1.Range.Text
$parenIdx = $txt.IndexOf("(")
Split-RunAt ($p1.Range.Start + $parenIdx) "SplitMark1"

# "six page" -> "one" + " page" (kept inside the existing gramStart/gramEnd
# proofErr markers, but as two separate runs).
$rngLimit = $d.Content
$rngLimit.Find.Execute("six page", $false, $false, $false, $false, $false, $true, 1, $false, "one page", 2)

$txt = $p1.Range.Text
$oneIdx = $txt.IndexOf("one page") + 3
Split-RunAt ($p1.Range.Start + $oneIdx) "SplitMark2"

# ===========================================================================
# 2) Remove the "_GoBack" bookmark paragraph and the "Research Experience"
#    Heading2 paragraph that followed the title -- the empty paragraph after
#    them is kept as-is.
# ===========================================================================
$pBookmark = $d.Paragraphs.Item(2)
$pResearchHeading = $d.Paragraphs.Item(3)
$killRange = $d.Range($pBookmark.Range.Start, $pResearchHeading.Range.End)
$killRange.Delete()

# Collapse the now-empty paragraph left behind into the following (originally
# empty) paragraph by removing its paragraph mark.
$pEmptyNow = $d.Paragraphs.Item(2)
$markRange = $d.Range($pEmptyNow.Range.End - 1, $pEmptyNow.Range.End)
$markRange.Delete()

# ===========================================================================
# 3) Remove every paragraph after the "Briefly summarize..." paragraph
#    through the end of the document (all the old guidance text/headings).
# ===========================================================================
$pAfterContent = $d.Paragraphs.Item(4)
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$tailRange = $d.Range($pAfterContent.Range.Start, $pLast.Range.End)
$tailRange.Delete()

# ===========================================================================
# 4) Replace the "Briefly summarize..." paragraph's text with the new
#    "Respective Contributions" guidance, formatted at 10pt (sz=20), with the
#    "_GoBack" bookmark reinstated in the middle of the text (this naturally
#    keeps the text split into the two runs the target markup expects).
# ===========================================================================
$pContent = $d.Paragraphs.Item(3)
$newText = "Describe the collaborative process between you and your sponsor/co-sponsor(s) in the development, review, and editing of this Research Training Plan. Also discuss your respective roles in accomplishing the proposed research."
$pContent.Range.Text = $newText

# Re-fetch the paragraph/range after the text swap and apply the 10pt size to
# the text only (excluding the trailing paragraph mark, so no pPr/rPr gets
# stamped on the paragraph mark itself).
$pContent = $d.Paragraphs.Item(3)
$textOnly = $d.Range($pContent.Range.Start, $pContent.Range.End - 1)
$textOnly.Font.Size = 10

# Insert the "_GoBack" bookmark right after "...Also di" (splitting "discuss"
# into "di" + "scuss"), matching the target markup.
$full = $pContent.Range.Text
$splitIdx = $full.IndexOf("Also di") + "Also di".Length
$goBackPos = $pContent.Range.Start + $splitIdx
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "[$i] $($p.Range.Text)"
}
